# Data Science Bootcamp Project #1 Plan - add a "Status" column replacing
# "Est. Completion", populate the first few status values, tidy up the
# column-A formatting quirk, fix the selection, and set the page to
# portrait orientation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Turn the "Est. Completion" header into "Status" and record the
#        status of the first few activities ------------------------------
$ws.Range("F1").Value = "Status"
$ws.Range("F2").Value = "Complete"
$ws.Range("F3").Value = "Complete"
$ws.Range("F4").Value = "In Progress"

# --- 2. Column A carried a stray column-level style (left over from
#        earlier formatting); clear it but keep every cell's own look
#        (bold section headers vs. indented task rows) intact -----------
$boldRows = @(1, 2, 6, 11, 17, 22)
$indentRows = @(3, 4, 5, 7, 8, 9, 10, 12, 13, 14, 15, 16, 18, 19, 20, 21, 23, 24, 25)

$ws.Columns.Item(1).ClearFormats()

foreach ($r in $boldRows) {
    $ws.Cells.Item($r, 1).Font.Bold = $true
}
foreach ($r in $indentRows) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.HorizontalAlignment = -4131
    $cell.IndentLevel = 1
}

# --- 3. Update the active selection to the "Understanding customer
#        behavior" block ---------------------------------------------------
$ws.Range("A11:A16").Select()

# --- 4. Set the sheet to print in portrait orientation --------------------
$ws.PageSetup.Orientation = 1

Write-Output "Applied Status column + formatting cleanup"
